# Update cryptocurrency price/volume data per latest scrape (Sat Dec 24 05:24:14 UTC 2022)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.55"
$ws.Range("D3").Value = "'21.99"
$ws.Range("D5").Value = "'0.05965"
$ws.Range("D6").Value = "'3.399"
$ws.Range("D7").Value = "'6.391"
$ws.Range("D9").Value = "'0.9658"
$ws.Range("D11").Value = "'0.03602"
$ws.Range("D12").Value = "'0.07406"
$ws.Range("D14").Value = "'0.09408"
$ws.Range("D15").Value = "'4.000"
$ws.Range("D16").Value = "'0.001596"
$ws.Range("D17").Value = "'0.04802"
$ws.Range("D18").Value = "'0.01121"
$ws.Range("E18").Value = "17OneONEBestin24h"
$ws.Range("D19").Value = "'0.006188"
$ws.Range("D20").Value = "'0.004141"
$ws.Range("D21").Value = "'0.0009834"
$ws.Range("D22").Value = "'0.00009705"
$ws.Range("D23").Value = "'3.743"
$ws.Range("D24").Value = "'2.165"
$ws.Range("D40").Value = "'0.03928"
$ws.Range("D41").Value = "'0.006514"
$ws.Range("D42").Value = "'0.1071"
$ws.Range("D43").Value = "'0.003002"
$ws.Range("D44").Value = "'0.005377"
$ws.Range("D45").Value = "'0.00005324"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.8505"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("D48").Value = "'0.04098"
